# This script applies a re-shuffle of the weekly price records in the
# "Fruta, Vega Monumental Concepción - Membrillo" sheet. Rows 2-32 each
# hold one record (columns D, L, M, N, O, P, Q, R, S, T). The edit
# rearranges which record occupies which row (a permutation of the 31
# rows) while columns A, B, C, E, F, G, H, I, J, K stay identical since
# they are constant across every row of this subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Destination row -> source row (i.e. after the edit, row $dest holds
# the data that used to live in row $map[$dest]).
$map = @{}
$map[2] = 19
$map[3] = 32
$map[4] = 31
$map[5] = 4
$map[6] = 30
$map[7] = 9
$map[8] = 23
$map[9] = 24
$map[10] = 27
$map[11] = 28
$map[12] = 13
$map[13] = 11
$map[14] = 12
$map[15] = 3
$map[16] = 17
$map[17] = 21
$map[18] = 6
$map[19] = 10
$map[20] = 2
$map[21] = 16
$map[22] = 25
$map[23] = 26
$map[24] = 5
$map[25] = 22
$map[26] = 15
$map[27] = 7
$map[28] = 8
$map[29] = 18
$map[30] = 20
$map[31] = 14
$map[32] = 29

# Columns that belong to each record (1-based column indices):
# D=4, L=12, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# 1) Snapshot the current value of every relevant cell before writing
#    anything, so overlapping/cyclical moves don't clobber data we
#    still need to read.
$buffer = @{}
for ($r = 2; $r -le 32; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $buffer[$r] = $rowData
}

# 2) Write the snapshotted values back out according to the mapping.
for ($dest = 2; $dest -le 32; $dest++) {
    $src = $map[$dest]
    $rowData = $buffer[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($dest, $c).Value = $rowData[$c]
    }
}
